$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse each card's multi-row fields (name + attribute rows) into a single
# Python-tuple-repr string per card, keeping only rows 1-6.

$ws.Range("A2").Value = "('Arena', ['Land', '{3}, {T}: Tap target creature you control and target creature of an opponent’s choice they control. Those creatures fight each other. (Each deals damage equal to its power to the other.)'])"
$ws.Range("A3").Value = "('Giant Badger', ['{1}{G}{G}', 'Creature — Badger', 'Whenever Giant Badger blocks, it gets +2/+2 until end of turn.', '2/2'])"
$ws.Range("A4").Value = "('Mana Crypt', ['{0}', 'Artifact', 'At the beginning of your upkeep, flip a coin. If you lose the flip, Mana Crypt deals 3 damage to you.', '{T}: Add {C}{C}.'])"
$ws.Range("A5").Value = "('Sewers of Estark', ['{2}{B}{B}', 'Instant', 'Choose target creature. If it’s attacking, it can’t be blocked this turn. If it’s blocking, prevent all combat damage that would be dealt this combat by it and each creature it’s blocking.'])"
$ws.Range("A6").Value = "('Windseeker Centaur', ['{1}{R}{R}', 'Creature — Centaur', 'Vigilance', '2/2'])"

# Remove the now-redundant rows that used to hold the individual fields.
$ws.Rows("7:23").Delete()
